$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 7.5
$ws.Range("I8").Value = 7.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 22.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 116.5
$ws.Range("N8").ClearContents()
$ws.Range("H18").Value = 723
$ws.Range("I18").Value = 723
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 723
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -439
$ws.Range("H70").Value = 3623.625
$ws.Range("I70").Value = 3500
$ws.Range("J70").Value = 3829.6667
$ws.Range("K70").Value = 10500
$ws.Range("L70").Value = 11489.0001
$ws.Range("M70").Value = -10230
$ws.Range("N70").Value = -12029.0001
$ws.Range("H73").Value = 3623.625
$ws.Range("I73").Value = 3500
$ws.Range("J73").Value = 3829.6667
$ws.Range("K73").Value = 10500
$ws.Range("L73").Value = 11489.0001
$ws.Range("M73").Value = -9564
$ws.Range("N73").Value = -13361.0001
$ws.Range("H80").Value = 2666.6667
$ws.Range("I80").Value = 2750
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 8250
$ws.Range("L80").Value = 7500
$ws.Range("M80").Value = -7252
$ws.Range("H83").Value = 2666.6667
$ws.Range("I83").Value = 2750
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 24750
$ws.Range("L83").Value = 22500
$ws.Range("M83").Value = -19758
$ws.Range("H137").Value = 5002
$ws.Range("I137").Value = 5002
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 15006
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -12456
$ws.Range("N137").ClearContents()
$ws.Range("H141").Value = 199
$ws.Range("I141").Value = 199
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 597
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 4583

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5992.5713
$ws.Range("I32").Value = 2536
$ws.Range("J32").Value = 18666.666
$ws.Range("K32").Value = 2536
$ws.Range("L32").Value = 18666.666
$ws.Range("M32").Value = -2249
$ws.Range("H61").Value = 3666.3333
$ws.Range("I61").Value = 3499.5
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3499.5
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -3287.5
$ws.Range("N61").Value = -4424
$ws.Range("H80").Value = 75164.664
$ws.Range("I80").Value = 79501
$ws.Range("J80").Value = 72996.5
$ws.Range("K80").Value = 79501
$ws.Range("L80").Value = 72996.5
$ws.Range("M80").Value = -78503
$ws.Range("N80").Value = -74992.5
$ws.Range("H83").Value = 75164.664
$ws.Range("I83").Value = 79501
$ws.Range("J83").Value = 72996.5
$ws.Range("K83").Value = 238503
$ws.Range("L83").Value = 218989.5
$ws.Range("M83").Value = -233511
$ws.Range("N83").Value = -228973.5
$ws.Range("H92").Value = 37800
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 37800
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 37800
$ws.Range("N92").Value = -42792
$ws.Range("H102").Value = 251572.25
$ws.Range("I102").Value = 500252
$ws.Range("J102").Value = 2892.5
$ws.Range("K102").Value = 500252
$ws.Range("L102").Value = 2892.5
$ws.Range("M102").Value = -498630
$ws.Range("H132").Value = 4228.778
$ws.Range("I132").Value = 4282
$ws.Range("J132").Value = 4042.5
$ws.Range("K132").Value = 12846
$ws.Range("L132").Value = 12127.5
$ws.Range("M132").Value = -10316
$ws.Range("N132").Value = -17187.5
$ws.Range("H136").Value = 3666.3333
$ws.Range("I136").Value = 3499.5
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 10498.5
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -7948.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 317.55554
$ws.Range("I11").Value = 75
$ws.Range("J11").Value = 438.83334
$ws.Range("K11").Value = 75
$ws.Range("L11").Value = 438.83334
$ws.Range("M11").Value = 65
$ws.Range("N11").Value = -718.83334
$ws.Range("H14").Value = 1999
$ws.Range("I14").Value = 1999
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1999
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1827
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H94").Value = 2817.1667
$ws.Range("I94").Value = 2880.9
$ws.Range("J94").Value = 2737.5
$ws.Range("K94").Value = 2880.9
$ws.Range("L94").Value = 2737.5
$ws.Range("M94").Value = -2429.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 9711.200000000001
$ws.Range("I11").Value = 1005
$ws.Range("J11").Value = 11887.75
$ws.Range("K11").Value = 1005
$ws.Range("L11").Value = 11887.75
$ws.Range("M11").Value = -865
$ws.Range("N11").Value = -12167.75
$ws.Range("H99").Value = 1002329.9
$ws.Range("I99").Value = 1001299.8
$ws.Range("J99").Value = 1003360
$ws.Range("K99").Value = 1001299.8
$ws.Range("L99").Value = 1003360
$ws.Range("M99").Value = -999801.8
$ws.Range("H126").Value = 1002329.9
$ws.Range("I126").Value = 1001299.8
$ws.Range("J126").Value = 1003360
$ws.Range("K126").Value = 3003899.4
$ws.Range("L126").Value = 3010080
$ws.Range("M126").Value = -3001429.4
$ws.Range("H132").Value = 3746.5881
$ws.Range("I132").Value = 2871.9285
$ws.Range("J132").Value = 7828.3335
$ws.Range("K132").Value = 8615.7855
$ws.Range("L132").Value = 23485.0005
$ws.Range("M132").Value = -6085.7855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 491.6
$ws.Range("I18").Value = 493
$ws.Range("J18").Value = 486
$ws.Range("K18").Value = 1479
$ws.Range("L18").Value = 1458
$ws.Range("M18").Value = -1310
$ws.Range("N18").Value = -1796
$ws.Range("H109").Value = 592.3333
$ws.Range("I109").Value = 592.3333
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 1776.9999
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -736.9999
$ws.Range("N109").ClearContents()
$ws.Range("H120").Value = 10833.333
$ws.Range("I120").Value = 3750
$ws.Range("J120").Value = 25000
$ws.Range("K120").Value = 11250
$ws.Range("L120").Value = 75000
$ws.Range("M120").Value = -6412
$ws.Range("N120").Value = -84676
$ws.Range("H129").Value = 1504.8889
$ws.Range("I129").Value = 380
$ws.Range("J129").Value = 1826.2858
$ws.Range("K129").Value = 1140
$ws.Range("L129").Value = 5478.857400000001
$ws.Range("M129").Value = 3860
$ws.Range("N129").Value = -15478.8574
$ws.Range("H137").Value = 2278.8
$ws.Range("I137").Value = 1000
$ws.Range("J137").Value = 2598.5
$ws.Range("K137").Value = 3000
$ws.Range("L137").Value = 7795.5
$ws.Range("M137").Value = 2100
$ws.Range("N137").Value = -17995.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1972
$ws.Range("I5").Value = 1972
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1972
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1860
$ws.Range("H92").Value = 21392.666
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 21392.666
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 21392.666
$ws.Range("N92").Value = -25136.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 495.75
$ws.Range("I46").Value = 428
$ws.Range("J46").Value = 699
$ws.Range("K46").Value = 428
$ws.Range("L46").Value = 699
$ws.Range("M46").Value = -240
$ws.Range("N46").Value = -1075
$ws.Range("H56").Value = 46000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 46000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 46000
$ws.Range("N56").Value = -47382
$ws.Range("M56").ClearContents()
$ws.Range("H109").Value = 42000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 42000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 42000
$ws.Range("N109").Value = -44774
$ws.Range("H136").Value = 638374.6
$ws.Range("I136").Value = 638374.6
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 1915123.8
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1912573.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 251525.5
$ws.Range("I2").Value = 334034
$ws.Range("J2").Value = 4000
$ws.Range("K2").Value = 334034
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = -333922
$ws.Range("N2").Value = -4224
$ws.Range("H4").Value = 2000.5
$ws.Range("I4").Value = 1001
$ws.Range("J4").Value = 3000
$ws.Range("K4").Value = 1001
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = -888
$ws.Range("N4").Value = -3226
$ws.Range("H64").Value = 10526
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 10526
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 10526
$ws.Range("N64").Value = -11022
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 10526
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 10526
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 10526
$ws.Range("N67").Value = -12242
$ws.Range("M67").ClearContents()
$ws.Range("H136").Value = 1470.7894
$ws.Range("I136").Value = 1290.8823
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3872.6469
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -1322.6469
